$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.432.94"
$ws.Range("E2").Value = "  +7.78%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.577.12"
$ws.Range("E3").Value = "  +9.77%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.04"
$ws.Range("E5").Value = "  +6.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.60"
$ws.Range("E6").Value = "  +7.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  +1.36%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.571.06"
$ws.Range("E9").Value = "  +9.36%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.10"
$ws.Range("E10").Value = "  +12.31%  "

# Row 11
$ws.Range("E11").Value = "  +6.24%  "

# Row 12
$ws.Range("E12").Value = "  +5.32%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("E13").Value = "  +1.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.019.64"
$ws.Range("E14").Value = "  +9.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.446.94"
$ws.Range("E15").Value = "  +7.91%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.81"
$ws.Range("E16").Value = "  +8.85%  "

# Row 17
$ws.Range("E17").Value = "  +5.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.580.37"
$ws.Range("E18").Value = "  +9.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  +3.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.49"
$ws.Range("E20").Value = "  +7.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.35"
$ws.Range("E21").Value = "  +7.65%  "

# Row 22
$ws.Range("E22").Value = "  +7.35%  "

# Row 23
$ws.Range("E23").Value = "  +0.35%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.92"
$ws.Range("E24").Value = "  +5.92%  "

# Row 25
$ws.Range("E25").Value = "  +5.61%  "

# Row 26
$ws.Range("E26").Value = "  +7.53%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.01"
$ws.Range("E27").Value = "  +0.95%  "

# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.685.37"
$ws.Range("E28").Value = "  +10.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0834"
$ws.Range("E29").Value = "  +11.30%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("E30").Value = "  +3.75%  "

# Row 31
$ws.Range("E31").Value = "  +0.16%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.31"
$ws.Range("E32").Value = "  +8.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.41"
$ws.Range("E33").Value = "  +6.66%  "

# Row 34
$ws.Range("E34").Value = "  +6.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.52"
$ws.Range("E35").Value = "  +8.05%  "

# Row 36
$ws.Range("E36").Value = "  +9.84%  "

# Row 37
$ws.Range("E37").Value = "  +8.90%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.853"
$ws.Range("E38").Value = "  +4.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.73"
$ws.Range("E39").Value = "  +9.86%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "298.98"
$ws.Range("E40").Value = "  +19.05%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.44"
$ws.Range("E41").Value = "  +7.51%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.19"
$ws.Range("E42").Value = "  +4.37%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0577"
$ws.Range("E43").Value = "  +10.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.101"
$ws.Range("E44").Value = "  +2.25%  "

# Row 45
$ws.Range("E45").Value = "  +9.11%  "

# Row 46
$ws.Range("E46").Value = "  +25.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").Value = "  -0.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.92"
$ws.Range("E48").Value = "  +12.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.21"
$ws.Range("E49").Value = "  +14.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("E50").Value = "  +6.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.25"
$ws.Range("E51").Value = "  +0.76%  "
